# Fruta / hortaliza, semanal
# Insert two new weekly-report rows (new row 33 and row 34) into the daily
# price log, pushing the existing rows 33..92 down to rows 35..94.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 33-92 down by inserting two blank rows at row 33.
$ws.Rows("33:34").Insert()

# --- New row 33 ---
$ws.Cells.Item(33, 1).Value = 1
$ws.Cells.Item(33, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(33, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(33, 4).Value = 44638
$ws.Cells.Item(33, 5).Value = 15
$ws.Cells.Item(33, 6).Value = "Fruta"
$ws.Cells.Item(33, 7).Value = 100109
$ws.Cells.Item(33, 8).Value = "Uva"
$ws.Cells.Item(33, 9).Value = 100109001
$ws.Cells.Item(33, 10).Value = "Uva"
$ws.Cells.Item(33, 11).Value = "Red Globe"
$ws.Cells.Item(33, 12).Value = "Segunda"
$ws.Cells.Item(33, 13).Value = 250
$ws.Cells.Item(33, 14).Value = 19000
$ws.Cells.Item(33, 15).Value = 20000
$ws.Cells.Item(33, 16).Value = 19500
$ws.Cells.Item(33, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(33, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(33, 19).Value = 975
$ws.Cells.Item(33, 20).Value = 20

# --- New row 34 ---
$ws.Cells.Item(34, 1).Value = 1
$ws.Cells.Item(34, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(34, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(34, 4).Value = 44638
$ws.Cells.Item(34, 5).Value = 15
$ws.Cells.Item(34, 6).Value = "Fruta"
$ws.Cells.Item(34, 7).Value = 100109
$ws.Cells.Item(34, 8).Value = "Uva"
$ws.Cells.Item(34, 9).Value = 100109001
$ws.Cells.Item(34, 10).Value = "Uva"
$ws.Cells.Item(34, 11).Value = "Superior Seedless"
$ws.Cells.Item(34, 12).Value = "Segunda"
$ws.Cells.Item(34, 13).Value = 270
$ws.Cells.Item(34, 14).Value = 19000
$ws.Cells.Item(34, 15).Value = 20000
$ws.Cells.Item(34, 16).Value = 19500
$ws.Cells.Item(34, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(34, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(34, 19).Value = 975
$ws.Cells.Item(34, 20).Value = 20
